$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Job")

# --- Header row updates (row 2 = type row, row 3 = table header row) ---
# Column O ("InitialLocked" bool) becomes "LevelNeed" (int)
$ws.Range("O2").Value = "int"
$ws.Range("O3").Value = "LevelNeed"

# --- Data rows: column O switches from a boolean "locked" flag to a numeric level requirement ---
$numberCells = @("O4","O6","O7","O9","O11","O12","O13","O14","O15","O16","O17")
foreach ($cellRef in $numberCells) {
    $ws.Range($cellRef).NumberFormat = "General"
}

$ws.Range("O4").ClearContents()
$ws.Range("O6").Value = 7
$ws.Range("O7").Value = 13
$ws.Range("O9").Value = 9
$ws.Range("O11").Value = 11
$ws.Range("O12").Value = 15
$ws.Range("O13").Value = 17
$ws.Range("O14").Value = 20
$ws.Range("O15").ClearContents()
$ws.Range("O16").ClearContents()
$ws.Range("O17").ClearContents()

# --- Column width tweak for the now-wider numeric column ---
$ws.Columns.Item(15).ColumnWidth = 6.375

# --- Move the active selection to O5, where the new "auto open job on level up" data lives ---
$ws.Range("O5").Select()
